# Update the cryptos list with refreshed prices / volume percentages.
# (Updated cryptos list on Tue Oct 24 07:49:00 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D / E store plain text (inline strings) in the source sheet, even when
# the text looks like a plain number (e.g. "0.999"). Assigning such a string to
# Range.Value would normally let Excel auto-convert it to a numeric cell, which
# would both change the cell type and can silently drop significant trailing
# zeros (e.g. "0.0520" -> 0.052). Prefixing with a leading apostrophe forces
# Excel to keep/store the value as text, exactly like the source data.
function Set-TextCell($address, $value) {
    $looksNumeric = $value -match '^[+-]?[0-9]*\.?[0-9]+$'
    if ($looksNumeric) {
        $ws.Range($address).Value = "'" + $value
    } else {
        $ws.Range($address).Value = $value
    }
}

# Row 2
Set-TextCell "D2" "33.858.39"
Set-TextCell "E2" "  +10.23%  "
# Row 3
Set-TextCell "D3" "1.804.35"
Set-TextCell "E3" "  +7.06%  "
# Row 4
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.09%  "
# Row 5
Set-TextCell "D5" "227.13"
Set-TextCell "E5" "  +3.04%  "
# Row 6
Set-TextCell "E6" "  +1.67%  "
# Row 7
Set-TextCell "D7" "0.999"
Set-TextCell "E7" "  -0.07%  "
# Row 8
Set-TextCell "D8" "31.03"
Set-TextCell "E8" "  +1.94%  "
# Row 9
Set-TextCell "D9" "47.10"
Set-TextCell "E9" "  +6.07%  "
# Row 10
Set-TextCell "D10" "0.279"
Set-TextCell "E10" "  +5.32%  "
# Row 11
Set-TextCell "E11" "  +5.48%  "
# Row 12
Set-TextCell "D12" "0.0927"
Set-TextCell "E12" "  +2.24%  "
# Row 13
Set-TextCell "D13" "2.062.51"
Set-TextCell "E13" "  +6.97%  "
# Row 14
Set-TextCell "D14" "1.808.30"
Set-TextCell "E14" "  +7.24%  "
# Row 15
Set-TextCell "E15" "  +1.80%  "
# Row 16
Set-TextCell "D16" "33.758.85"
Set-TextCell "E16" "  +9.88%  "
# Row 17
Set-TextCell "D17" "10.04"
Set-TextCell "E17" "  -3.85%  "
# Row 18
Set-TextCell "D18" "4.23"
Set-TextCell "E18" "  +6.06%  "
# Row 19
Set-TextCell "D19" "68.99"
Set-TextCell "E19" "  +3.87%  "
# Row 20
Set-TextCell "D20" "254.78"
Set-TextCell "E20" "  +3.61%  "
# Row 21
Set-TextCell "D21" "0.0₃0739"
Set-TextCell "E21" "  +3.30%  "
# Row 22
Set-TextCell "E22" "  +0.14%  "
# Row 23
Set-TextCell "D23" "10.39"
Set-TextCell "E23" "  +1.84%  "
# Row 24
Set-TextCell "E24" "  +0.10%  "
# Row 25
Set-TextCell "E25" "  +0.53%  "
# Row 26
Set-TextCell "D26" "156.52"
Set-TextCell "E26" "  -1.01%  "
# Row 27
Set-TextCell "B27" "MinaProtocolToken"
Set-TextCell "C27" "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
Set-TextCell "D27" "2.76"
Set-TextCell "E27" "  +575.43%  "
# Row 28
Set-TextCell "B28" "EthereumClassic"
Set-TextCell "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D28" "16.37"
Set-TextCell "E28" "  +3.04%  "
# Row 29
Set-TextCell "B29" "Stellar"
Set-TextCell "C29" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D29" "0.114"
Set-TextCell "E29" "  +2.51%  "
# Row 30
Set-TextCell "B30" "Cosmos"
Set-TextCell "C30" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D30" "7.01"
Set-TextCell "E30" "  +4.62%  "
# Row 31
Set-TextCell "D31" "0.999"
Set-TextCell "E31" "  -0.07%  "
# Row 32
Set-TextCell "D32" "3.81"
Set-TextCell "E32" "  +9.33%  "
# Row 33
Set-TextCell "D33" "0.0506"
Set-TextCell "E33" "  +1.50%  "
# Row 34
Set-TextCell "D34" "1.20"
Set-TextCell "E34" "  +4.82%  "
# Row 35
Set-TextCell "D35" "3.47"
Set-TextCell "E35" "  +5.62%  "
# Row 36
Set-TextCell "D36" "1.531.39"
Set-TextCell "E36" "  +1.16%  "
# Row 37
Set-TextCell "E37" "  +2.03%  "
# Row 38
Set-TextCell "D38" "1.07"
Set-TextCell "E38" "  +3.47%  "
# Row 39
Set-TextCell "D39" "83.32"
Set-TextCell "E39" "  -1.52%  "
# Row 40
Set-TextCell "E40" "  +3.45%  "
# Row 41
Set-TextCell "D41" "0.611"
Set-TextCell "E41" "  +3.99%  "
# Row 42
Set-TextCell "E42" "  +2.79%  "
# Row 43
Set-TextCell "D43" "2.32"
Set-TextCell "E43" "  +0.98%  "
# Row 44
Set-TextCell "D44" "0.901"
Set-TextCell "E44" "  +7.39%  "
# Row 45
Set-TextCell "D45" "2.11"
Set-TextCell "E45" "  +5.85%  "
# Row 46
Set-TextCell "D46" "0.0520"
Set-TextCell "E46" "  +4.05%  "
# Row 47
Set-TextCell "E47" "  +3.99%  "
# Row 48
Set-TextCell "D48" "1.952.12"
Set-TextCell "E48" "  +6.97%  "
# Row 49
Set-TextCell "D49" "0.998"
Set-TextCell "E49" "  -0.11%  "
# Row 50
Set-TextCell "D50" "5.64"
Set-TextCell "E50" "  +3.22%  "
# Row 51
Set-TextCell "D51" "51.93"
Set-TextCell "E51" "  -0.01%  "
